# From v1.0.2 to v1.0.3
# Reorders shared test-case step content:
#  - TC3 and TC4 step contents are swapped (TC3 becomes "realizar o empenho",
#    TC4 becomes "atribuir/desatribuir")
#  - The "Expected Results" text on TC1's second step (D11) and TC6's step (D48)
#    are swapped (success message moves to D11, the MSG207 error message moves to D48)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap D11 <-> D48
$ws.Range("D11").Value = "SYSTEM Realiza a persistência do número do credor, para o beneficiário indicado, na base do RH; Atualiza a listagem de solicitações aguardando empenho, já com o número do credor recém informado; e Exibe mensagem de sucesso."
$ws.Range("D48").Value = "SYSTEM Identifica que houve um erro inesperado, quando da tentativa de inserção do número do credor; e Exibe mensagem de erro (MSG207 - Não foi possível atualizar o número do credor) para o usuário."

# Swap TC3 (row 26) <-> TC4 (row 33) step contents
$ws.Range("B26").Value = "Chefe/Beneficiário Clica para realizar o empenho de uma diária."
$ws.Range("D26").Value = "SYSTEM Apresenta a tela de Registrar Empenho"

$ws.Range("B33").Value = "Chefe/Beneficiário Clica para atribuir/desatribuir o registro a si mesmo."
$ws.Range("D33").Value = "SYSTEM Atualiza a lista de registros de solicitações, onde o nome deverá constar o nome do usuário logado (que se atribuiu como responsável pelo empenho), no campo de atribuição (no caso de desatribuição, o nome deverá ser removido)."
